# Export latest mex files to xlsx. Update segments and articles_db.
# Appends 4 new coded-segment rows (213-216) to Sheet1, mirroring the
# formatting of the existing data rows (row 212 is used as the template).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10
$dot = [char]0x25CF

# Clone the formatting (styles/borders/fills) of the last existing data
# row (212) onto the four new rows so every column keeps its established
# look (colored bullet, text columns, integer/percent columns, etc).
$ws.Range("A212:M212").Copy()
$ws.Range("A213:M216").PasteSpecial(-4122)

# ---------------------------------------------------------------- Row 213
$ws.Range("A213").Value = $dot
$ws.Range("D213").Value = "'3095"
$ws.Range("D212").Copy()
$ws.Range("D213").PasteSpecial(-4122)
$ws.Range("E213").Value = "Bacteria:Binomial (genus species)"
$ws.Range("F213").Value = "1: 1498"
$ws.Range("G213").Value = "1: 1519"
$ws.Range("H213").Value = 0
$ws.Range("I213").Value = "Pseudomonas aeruginosa"
$ws.Range("J213").Value = 22
$ws.Range("K213").Value = 0.14659800000000001
$ws.Range("L213").Value = "Sonia"
$ws.Range("M213").Value = "11/8/18 14:15:00"
$ws.Rows.Item(213).RowHeight = 16

# ---------------------------------------------------------------- Row 214
$ws.Range("A214").Value = $dot
$ws.Range("D214").Value = "'3557"
$ws.Range("D212").Copy()
$ws.Range("D214").PasteSpecial(-4122)
$ws.Range("E214").Value = "Bacteria:Binomial (genus species)"
$ws.Range("F214").Value = "1: 44"
$ws.Range("G214").Value = "1: 61"
$ws.Range("H214").Value = 0
$ws.Range("I214").Value = "Enterobacteriaceae"
$ws.Range("J214").Value = 18
$ws.Range("K214").Value = 0.057685
$ws.Range("L214").Value = "Sonia"
$ws.Range("M214").Value = "11/8/18 14:18:00"
$ws.Rows.Item(214).RowHeight = 16

# ---------------------------------------------------------------- Row 215
$ws.Range("A215").Value = $dot
$ws.Range("D215").Value = "'18527"
$ws.Range("D212").Copy()
$ws.Range("D215").PasteSpecial(-4122)
$ws.Range("E215").Value = "Bacteria:Binomial (genus species)"
$ws.Range("F215").Value = "1: 2966"
$ws.Range("G215").Value = "1: 2990"
$ws.Range("H215").Value = 0
$ws.Range("I215").Value = "cinetobacter  " + $nl + "baumannii "
$ws.Range("J215").Value = 23
$ws.Range("K215").Value = 0.138881
$ws.Range("L215").Value = "Sonia"
$ws.Range("M215").Value = "11/8/18 14:18:00"
$ws.Rows.Item(215).RowHeight = 30

# ---------------------------------------------------------------- Row 216
$ws.Range("A216").Value = $dot
$ws.Range("D216").Value = "'15760"
$ws.Range("D212").Copy()
$ws.Range("D216").PasteSpecial(-4122)
$ws.Range("E216").Value = "Bacteria:Binomial (genus species)"
$ws.Range("F216").Value = "1: 2676"
$ws.Range("G216").Value = "1: 2695"
$ws.Range("H216").Value = 0
$ws.Range("I216").Value = "Oligella ureolytica "
$ws.Range("J216").Value = 19
$ws.Range("K216").Value = 0.12753400000000001
$ws.Range("L216").Value = "Sonia"
$ws.Range("M216").Value = "11/8/18 14:20:00"
$ws.Rows.Item(216).RowHeight = 16
